# Auto commit at 2025-12-19  8:13:38.73
# Updates the "Metrics" sheet's source figures (B2:B13) for the latest
# refresh. Cells on "today" (and any other sheet) that reference these
# via formulas recalc automatically.

$wb = $excel.ActiveWorkbook

$wsMetrics = $wb.Worksheets.Item("Metrics")
$wsToday   = $wb.Worksheets.Item("today")

$wsMetrics.Range("B2").Value  = 241820.25
$wsMetrics.Range("B3").Value  = 207154.29
$wsMetrics.Range("B4").Value  = 74222.5
$wsMetrics.Range("B5").Value  = 9908
$wsMetrics.Range("B6").Value  = 5444527.3600000013
$wsMetrics.Range("B7").Value  = 4607507.25
$wsMetrics.Range("B8").Value  = 1606179.3800000004
$wsMetrics.Range("B9").Value  = 212615
$wsMetrics.Range("B10").Value = 33909908.349999987
$wsMetrics.Range("B11").Value = 31882782.41
$wsMetrics.Range("B12").Value = 11887901.419999996
$wsMetrics.Range("B13").Value = 1310245

# Move the cell selection on the "Metrics" sheet (it is not the active
# tab), then restore "today" as the active sheet/tab afterwards.
$wsMetrics.Range("D20").Select() | Out-Null
$wsToday.Activate() | Out-Null
